$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("E4")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.Font.Color = 0
$rng.Font.Family = 3
$rng.WrapText = $true
$rng.VerticalAlignment = -4108
$rng6 = $ws.Range("E6")

$rng6.Font.Name = "宋体"
$rng6.Font.Size = 12
$rng6.Font.Color = 0
$rng6.Font.Family = 3
$rng6.WrapText = $true
$rng6.VerticalAlignment = -4108